# Refresh the cryptos list: updated Price (col D) and Volume(1h) (col E) figures,
# plus two coin-row swaps (Stellar/WEMIXTOKEN and Algorand/VeChain) from the
# scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some prices are plain decimals (e.g. "0.9954", "1.000") which Excel's COM
# layer would otherwise auto-coerce into numbers, losing the trailing zeros /
# text formatting the source data relies on. Force those through as text via
# a temporary Text number format, then drop the format again so the cell
# keeps its original (default) style.
$ws.Range('D2').Value = '26.055.68'
$ws.Range('E2').Value = '  +6.25%  '

$ws.Range('D3').Value = '1.713.80'
$ws.Range('E3').Value = '  +3.74%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9954'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.61%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9984'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3685'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.79%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.43'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.95%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3322'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.15%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.178'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.58%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07494'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.20%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.224'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.11'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.895'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.27%  '

$ws.Range('D16').Value = '1.701.57'
$ws.Range('E16').Value = '  +3.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001075'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.10%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06637'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.52'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.75%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9974'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.23'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.37%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.059'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.11%  '

$ws.Range('E23').Value = '  +4.04%  '

$ws.Range('D24').Value = '25.922.87'
$ws.Range('E24').Value = '  +5.83%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.463'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.34%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.490'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.91%  '

$ws.Range('E28').Value = '  +3.50%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.309'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +10.19%  '

$ws.Range('D30').Value = '1.893.42'
$ws.Range('E30').Value = '  +3.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.40'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.36%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.108'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.964'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.36%  '

$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.721'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.82%  '

$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08509'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.76%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.92'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.349'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.75%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.276'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06211'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.05%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.544'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.46%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02277'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.86%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2121'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.66'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +16.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6169'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.58%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9982'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.837'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.29%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5859'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.10'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.006'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07255'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.86%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.71'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.25%  '
